$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Corrected life_expectancy values (column C), female rows (D = "Female") ---
$ws.Range("C2").Value  = 79.3    # age 0
$ws.Range("C3").Value  = 78.7    # age 1
$ws.Range("C5").Value  = 69.8    # age 10
$ws.Range("C17").Value = 16      # age 70
$ws.Range("C18").Value = 12.5    # age 75
$ws.Range("C19").Value = 9.4     # age 80
$ws.Range("C20").Value = 6.7     # age 85

# --- Corrected life_expectancy values (column C), male rows (D = "Male") ---
$ws.Range("C21").Value = 73.5    # age 0
$ws.Range("C22").Value = 73      # age 1
$ws.Range("C23").Value = 69.1    # age 5
$ws.Range("C24").Value = 64.1    # age 10
$ws.Range("C25").Value = 59.2    # age 15
$ws.Range("C26").Value = 54.4    # age 20
$ws.Range("C27").Value = 49.8    # age 25
$ws.Range("C28").Value = 45.4    # age 30
$ws.Range("C29").Value = 41      # age 35
$ws.Range("C30").Value = 36.6    # age 40
$ws.Range("C31").Value = 32.3    # age 45
$ws.Range("C32").Value = 28.2    # age 50
$ws.Range("C33").Value = 24.2    # age 55
$ws.Range("C38").Value = 7.9     # age 80

# --- View state: scroll down and reselect the cell below the last row of data ---
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("C40").Select()
